$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.144900666666667
$ws.Range("H2").Value = 3.434702
$ws.Range("I2").Value = 0.02523133726002265
$ws.Range("J2").Value = 0.02523133726002265
$ws.Range("M2").Value = 11.61289466666667
$ws.Range("N2").Value = 34.838684
$ws.Range("O2").Value = 0.09693042549509606
$ws.Range("P2").Value = 0.09693042549509606
$ws.Range("Q2").Value = 13.29561084579644
$ws.Range("R2").Value = 119.660497612168
$ws.Range("S2").Value = 0.002445684256424266
$ws.Range("T2").Value = 0.002445684256424266
$ws.Range("G3").Value = 1.144900666666667
$ws.Range("H3").Value = 3.434702
$ws.Range("I3").Value = 0.02523133726002265
$ws.Range("J3").Value = 0.02523133726002265
$ws.Range("O3").Value = 0.2981108740043866
$ws.Range("P3").Value = 0.2981108740043866
$ws.Range("Q3").Value = 40.89083638514622
$ws.Range("R3").Value = 368.017527466316
$ws.Range("S3").Value = 0.007521736002884795
$ws.Range("T3").Value = 0.007521736002884795
$ws.Range("G4").Value = 1.144900666666667
$ws.Range("H4").Value = 3.434702
$ws.Range("I4").Value = 0.02523133726002265
$ws.Range("J4").Value = 0.02523133726002265
$ws.Range("M4").Value = 27.39934733333333
$ws.Range("N4").Value = 82.198042
$ws.Range("O4").Value = 0.2286966748205465
$ws.Range("P4").Value = 0.2286966748205465
$ws.Range("Q4").Value = 31.36953102816489
$ws.Range("R4").Value = 282.325779253484
$ws.Range("S4").Value = 0.005770322932642938
$ws.Range("T4").Value = 0.005770322932642938
$ws.Range("G5").Value = 1.144900666666667
$ws.Range("H5").Value = 3.434702
$ws.Range("I5").Value = 0.02523133726002265
$ws.Range("J5").Value = 0.02523133726002265
$ws.Range("M5").Value = 45.078635
$ws.Range("N5").Value = 135.235905
$ws.Range("O5").Value = 0.3762620256799708
$ws.Range("P5").Value = 0.3762620256799709
$ws.Range("Q5").Value = 51.61055926392333
$ws.Range("R5").Value = 464.49503337531
$ws.Range("S5").Value = 0.009493594068070644
$ws.Range("T5").Value = 0.009493594068070646
$ws.Range("I6").Value = 0.03646539869776051
$ws.Range("J6").Value = 0.03646539869776051
$ws.Range("M6").Value = 11.61289466666667
$ws.Range("N6").Value = 34.838684
$ws.Range("O6").Value = 0.09693042549509606
$ws.Range("P6").Value = 0.09693042549509606
$ws.Range("Q6").Value = 19.215380676252
$ws.Range("R6").Value = 172.938426086268
$ws.Range("S6").Value = 0.003534606611622248
$ws.Range("T6").Value = 0.003534606611622248
$ws.Range("I7").Value = 0.03646539869776051
$ws.Range("J7").Value = 0.03646539869776051
$ws.Range("O7").Value = 0.2981108740043866
$ws.Range("P7").Value = 0.2981108740043866
$ws.Range("S7").Value = 0.01087073187670781
$ws.Range("T7").Value = 0.01087073187670781
$ws.Range("I8").Value = 0.03646539869776051
$ws.Range("J8").Value = 0.03646539869776051
$ws.Range("M8").Value = 27.39934733333333
$ws.Range("N8").Value = 82.198042
$ws.Range("O8").Value = 0.2286966748205465
$ws.Range("P8").Value = 0.2286966748205465
$ws.Range("Q8").Value = 45.336576659226
$ws.Range("R8").Value = 408.029189933034
$ws.Range("S8").Value = 0.008339515428183317
$ws.Range("T8").Value = 0.008339515428183317
$ws.Range("I9").Value = 0.03646539869776051
$ws.Range("J9").Value = 0.03646539869776051
$ws.Range("M9").Value = 45.078635
$ws.Range("N9").Value = 135.235905
$ws.Range("O9").Value = 0.3762620256799708
$ws.Range("P9").Value = 0.3762620256799709
$ws.Range("Q9").Value = 74.58976911046499
$ws.Range("R9").Value = 671.307921994185
$ws.Range("S9").Value = 0.01372054478124714
$ws.Range("T9").Value = 0.01372054478124714
$ws.Range("G10").Value = 3.191626333333333
$ws.Range("H10").Value = 9.574878999999999
$ws.Range("I10").Value = 0.07033710676294723
$ws.Range("J10").Value = 0.07033710676294723
$ws.Range("M10").Value = 11.61289466666667
$ws.Range("N10").Value = 34.838684
$ws.Range("O10").Value = 0.09693042549509606
$ws.Range("P10").Value = 0.09693042549509606
$ws.Range("Q10").Value = 37.06402042435955
$ws.Range("R10").Value = 333.576183819236
$ws.Range("S10").Value = 0.006817805686626474
$ws.Range("T10").Value = 0.006817805686626474
$ws.Range("G11").Value = 3.191626333333333
$ws.Range("H11").Value = 9.574878999999999
$ws.Range("I11").Value = 0.07033710676294723
$ws.Range("J11").Value = 0.07033710676294723
$ws.Range("O11").Value = 0.2981108740043866
$ws.Range("P11").Value = 0.2981108740043866
$ws.Range("Q11").Value = 113.9909111755758
$ws.Range("R11").Value = 1025.918200580182
$ws.Range("S11").Value = 0.02096825637204205
$ws.Range("T11").Value = 0.02096825637204205
$ws.Range("G12").Value = 3.191626333333333
$ws.Range("H12").Value = 9.574878999999999
$ws.Range("I12").Value = 0.07033710676294723
$ws.Range("J12").Value = 0.07033710676294723
$ws.Range("M12").Value = 27.39934733333333
$ws.Range("N12").Value = 82.198042
$ws.Range("O12").Value = 0.2286966748205465
$ws.Range("P12").Value = 0.2286966748205465
$ws.Range("Q12").Value = 87.44847846521311
$ws.Range("R12").Value = 787.036306186918
$ws.Range("S12").Value = 0.01608586243318381
$ws.Range("T12").Value = 0.01608586243318381
$ws.Range("G13").Value = 3.191626333333333
$ws.Range("H13").Value = 9.574878999999999
$ws.Range("I13").Value = 0.07033710676294723
$ws.Range("J13").Value = 0.07033710676294723
$ws.Range("M13").Value = 45.078635
$ws.Range("N13").Value = 135.235905
$ws.Range("O13").Value = 0.3762620256799708
$ws.Range("P13").Value = 0.3762620256799709
$ws.Range("Q13").Value = 143.8741585367217
$ws.Range("R13").Value = 1294.867426830495
$ws.Range("S13").Value = 0.0264651822710949
$ws.Range("T13").Value = 0.0264651822710949
$ws.Range("G14").Value = 39.384953
$ws.Range("H14").Value = 118.154859
$ws.Range("I14").Value = 0.8679661572792696
$ws.Range("J14").Value = 0.8679661572792696
$ws.Range("M14").Value = 11.61289466666667
$ws.Range("N14").Value = 34.838684
$ws.Range("O14").Value = 0.09693042549509606
$ws.Range("P14").Value = 0.09693042549509606
$ws.Range("Q14").Value = 457.3733106406173
$ws.Range("R14").Value = 4116.359795765556
$ws.Range("S14").Value = 0.08413232894042308
$ws.Range("T14").Value = 0.08413232894042308
$ws.Range("G15").Value = 39.384953
$ws.Range("H15").Value = 118.154859
$ws.Range("I15").Value = 0.8679661572792696
$ws.Range("J15").Value = 0.8679661572792696
$ws.Range("O15").Value = 0.2981108740043866
$ws.Range("P15").Value = 0.2981108740043866
$ws.Range("Q15").Value = 1406.657988809225
$ws.Range("R15").Value = 12659.92189928302
$ws.Range("S15").Value = 0.2587501497527519
$ws.Range("T15").Value = 0.2587501497527519
$ws.Range("G16").Value = 39.384953
$ws.Range("H16").Value = 118.154859
$ws.Range("I16").Value = 0.8679661572792696
$ws.Range("J16").Value = 0.8679661572792696
$ws.Range("M16").Value = 27.39934733333333
$ws.Range("N16").Value = 82.198042
$ws.Range("O16").Value = 0.2286966748205465
$ws.Range("P16").Value = 0.2286966748205465
$ws.Range("Q16").Value = 1079.122006954009
$ws.Range("R16").Value = 9712.098062586078
$ws.Range("S16").Value = 0.1985009740265365
$ws.Range("T16").Value = 0.1985009740265365
$ws.Range("G17").Value = 39.384953
$ws.Range("H17").Value = 118.154859
$ws.Range("I17").Value = 0.8679661572792696
$ws.Range("J17").Value = 0.8679661572792696
$ws.Range("M17").Value = 45.078635
$ws.Range("N17").Value = 135.235905
$ws.Range("O17").Value = 0.3762620256799708
$ws.Range("P17").Value = 0.3762620256799709
$ws.Range("Q17").Value = 1775.419920779155
$ws.Range("R17").Value = 15978.77928701239
$ws.Range("S17").Value = 0.3265827045595581
$ws.Range("T17").Value = 0.3265827045595582
